# adjust error code and message for connector test with not exist order input
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("getConceptModelDataByCondition")
$ws.Activate()

# Rows 17 and 18: rspMessage (column M) changed from
# "fields false not exist in entity!" to "not exist in entity!"
$ws.Range("M17").Value = "not exist in entity!"
$ws.Range("M18").Value = "not exist in entity!"

# Rows 23-26: rspMessage (column M) changed to the same new text, and
# rspCode (column L) changed from 106601 to 106107.
$ws.Range("M23").Value = "not exist in entity!"
$ws.Range("L23").Value = 106107

$ws.Range("M24").Value = "not exist in entity!"
$ws.Range("L24").Value = 106107

$ws.Range("M25").Value = "not exist in entity!"
$ws.Range("L25").Value = 106107

$ws.Range("M26").Value = "not exist in entity!"
$ws.Range("L26").Value = 106107

# Leave the cursor where the author last left it when saving.
$ws.Range("M20").Select()
